# The "work areas" list in column A contained a handful of labels using a
# Unicode EN DASH (U+2013) as a separator, e.g. "FUSELAGE – CABIN COMPARTMENT".
# That multi-byte character was not supported by the destination database's
# character set (utf8mb4 import), so every EN DASH in the sheet is normalized
# to a plain ASCII hyphen-minus ("-").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:A171")

# Replace the EN DASH character (Unicode 0x2013) with a regular hyphen
# everywhere it appears in the work-area names.
[void]$dataRange.Replace([char]0x2013, "-")

# The sheet is maintained as an alphabetically sorted list (see the existing
# <sortState> on the sheet, sorted on column A). Re-apply the same ascending
# sort so the renamed rows land back in their correct alphabetical position.
[void]$dataRange.Sort($ws.Range("A1"), 1)

# Leave the selection near where the data was last touched.
$ws.Range("K102").Select() | Out-Null
